# Updated symbol list (GitHub Actions crypto price refresh).
# Column D (Price) and column G (Hora) are plain text cells in the sheet
# (stored as inline strings, numeric-looking though they are), so numeric
# literals are written with a leading apostrophe to force text entry and
# avoid Excel silently re-typing them as numbers. Columns B/C/E are
# ordinary text, so no prefix is needed there.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.45"
$ws.Range("G2").Value = "'3"
$ws.Range("D3").Value = "'23.18"
$ws.Range("G3").Value = "'3"
$ws.Range("D4").Value = "'5.421"
$ws.Range("G4").Value = "'3"
$ws.Range("D5").Value = "'0.06046"
$ws.Range("G5").Value = "'3"
$ws.Range("D6").Value = "'3.403"
$ws.Range("G6").Value = "'3"
$ws.Range("D7").Value = "'0.8046"
$ws.Range("G7").Value = "'3"
$ws.Range("D8").Value = "'0.9322"
$ws.Range("G8").Value = "'3"
$ws.Range("D9").Value = "'0.1429"
$ws.Range("G9").Value = "'3"
$ws.Range("D10").Value = "'0.07463"
$ws.Range("G10").Value = "'3"
$ws.Range("D11").Value = "'0.03363"
$ws.Range("G11").Value = "'3"
$ws.Range("G12").Value = "'3"
$ws.Range("D13").Value = "'4.009"
$ws.Range("G13").Value = "'3"
$ws.Range("D14").Value = "'0.09365"
$ws.Range("G14").Value = "'3"
$ws.Range("G15").Value = "'3"
$ws.Range("D16").Value = "'0.04836"
$ws.Range("G16").Value = "'3"
$ws.Range("G17").Value = "'3"
$ws.Range("D18").Value = "'0.005100"
$ws.Range("G18").Value = "'3"
$ws.Range("D19").Value = "'0.004160"
$ws.Range("G19").Value = "'3"
$ws.Range("D20").Value = "'0.0009869"
$ws.Range("G20").Value = "'3"
$ws.Range("D21").Value = "'0.00008705"
$ws.Range("G21").Value = "'3"
$ws.Range("D22").Value = "'3.655"
$ws.Range("G22").Value = "'3"
$ws.Range("D23").Value = "'6.445"
$ws.Range("G23").Value = "'3"
$ws.Range("D24").Value = "'2.187"
$ws.Range("G24").Value = "'3"
$ws.Range("G25").Value = "'3"
$ws.Range("G26").Value = "'3"
$ws.Range("G27").Value = "'3"
$ws.Range("G28").Value = "'3"
$ws.Range("G29").Value = "'3"
$ws.Range("G30").Value = "'3"
$ws.Range("G31").Value = "'3"
$ws.Range("G32").Value = "'3"
$ws.Range("G33").Value = "'3"
$ws.Range("G34").Value = "'3"
$ws.Range("G35").Value = "'3"
$ws.Range("G36").Value = "'3"
$ws.Range("G37").Value = "'3"
$ws.Range("G38").Value = "'3"
$ws.Range("G39").Value = "'3"
$ws.Range("D40").Value = "'0.03979"
$ws.Range("G40").Value = "'3"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006424"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("G41").Value = "'3"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1074"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("G42").Value = "'3"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002902"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("G43").Value = "'3"
$ws.Range("D44").Value = "'0.006300"
$ws.Range("G44").Value = "'3"
$ws.Range("G45").Value = "'3"
$ws.Range("G46").Value = "'3"
$ws.Range("G47").Value = "'3"
$ws.Range("D48").Value = "'0.9005"
$ws.Range("E48").Value = "47CoinbaseStockTokenCOIN"
$ws.Range("G48").Value = "'3"
$ws.Range("G49").Value = "'3"
$ws.Range("G50").Value = "'3"
$ws.Range("D51").Value = "'0.01011"
$ws.Range("G51").Value = "'3"
